$wb = $excel.ActiveWorkbook

# Sheet 1: ROW50-FE-LIFTER (rows 29-31)
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(29,1).Value = 45729.58035023148
$ws.Cells.Item(29,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(29,2).Value = "0x01,0x90"
$ws.Cells.Item(29,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
$ws.Cells.Item(29,4).Value = "0x01,0x90,"
$ws.Cells.Item(29,5).Value = "0x14"
$ws.Cells.Item(29,6).Value = 400
$ws.Cells.Item(29,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(29,8).Value = 400
$ws.Cells.Item(29,9).Value = 20

$ws.Cells.Item(30,1).Value = 45729.5803721875
$ws.Cells.Item(30,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(30,2).Value = "0x01,0x90"
$ws.Cells.Item(30,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
$ws.Cells.Item(30,4).Value = "0x01,0x90,"
$ws.Cells.Item(30,5).Value = "0x14"
$ws.Cells.Item(30,6).Value = 400
$ws.Cells.Item(30,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(30,8).Value = 400
$ws.Cells.Item(30,9).Value = 20

$ws.Cells.Item(31,1).Value = 45729.58039546297
$ws.Cells.Item(31,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(31,2).Value = "0x01,0x90"
$ws.Cells.Item(31,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
$ws.Cells.Item(31,4).Value = "0x01,0x90,"
$ws.Cells.Item(31,5).Value = "0x14"
$ws.Cells.Item(31,6).Value = 400
$ws.Cells.Item(31,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(31,8).Value = 400
$ws.Cells.Item(31,9).Value = 20

# Sheet 2: ROW50-MID-LIFTER (rows 68-85)
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(68,1).Value = 45729.31518523148
$ws.Cells.Item(68,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(68,2).Value = "0x01,0x90"
$ws.Cells.Item(68,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(68,4).Value = "0x01,0x90,"
$ws.Cells.Item(68,5).Value = "0x19"
$ws.Cells.Item(68,6).Value = 400
$ws.Cells.Item(68,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(68,8).Value = 400
$ws.Cells.Item(68,9).Value = 25

$ws.Cells.Item(69,1).Value = 45729.31520722222
$ws.Cells.Item(69,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(69,2).Value = "0x01,0x90"
$ws.Cells.Item(69,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(69,4).Value = "0x01,0x90,"
$ws.Cells.Item(69,5).Value = "0x19"
$ws.Cells.Item(69,6).Value = 400
$ws.Cells.Item(69,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(69,8).Value = 400
$ws.Cells.Item(69,9).Value = 25

$ws.Cells.Item(70,1).Value = 45729.3152303588
$ws.Cells.Item(70,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(70,2).Value = "0x01,0x90"
$ws.Cells.Item(70,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(70,4).Value = "0x01,0x90,"
$ws.Cells.Item(70,5).Value = "0x19"
$ws.Cells.Item(70,6).Value = 400
$ws.Cells.Item(70,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(70,8).Value = 400
$ws.Cells.Item(70,9).Value = 25

$ws.Cells.Item(71,1).Value = 45729.39866047454
$ws.Cells.Item(71,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(71,2).Value = "0x01,0x90"
$ws.Cells.Item(71,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(71,4).Value = "0x01,0x90,"
$ws.Cells.Item(71,5).Value = "0x19"
$ws.Cells.Item(71,6).Value = 400
$ws.Cells.Item(71,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(71,8).Value = 400
$ws.Cells.Item(71,9).Value = 25

$ws.Cells.Item(72,1).Value = 45729.39868246527
$ws.Cells.Item(72,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(72,2).Value = "0x01,0x90"
$ws.Cells.Item(72,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(72,4).Value = "0x01,0x90,"
$ws.Cells.Item(72,5).Value = "0x19"
$ws.Cells.Item(72,6).Value = 400
$ws.Cells.Item(72,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(72,8).Value = 400
$ws.Cells.Item(72,9).Value = 25

$ws.Cells.Item(73,1).Value = 45729.39870561343
$ws.Cells.Item(73,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(73,2).Value = "0x01,0x90"
$ws.Cells.Item(73,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(73,4).Value = "0x01,0x90,"
$ws.Cells.Item(73,5).Value = "0x19"
$ws.Cells.Item(73,6).Value = 400
$ws.Cells.Item(73,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(73,8).Value = 400
$ws.Cells.Item(73,9).Value = 25

$ws.Cells.Item(74,1).Value = 45729.48213673611
$ws.Cells.Item(74,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(74,2).Value = "0x01,0x90"
$ws.Cells.Item(74,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(74,4).Value = "0x01,0x90,"
$ws.Cells.Item(74,5).Value = "0x19"
$ws.Cells.Item(74,6).Value = 400
$ws.Cells.Item(74,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(74,8).Value = 400
$ws.Cells.Item(74,9).Value = 25

$ws.Cells.Item(75,1).Value = 45729.48215891204
$ws.Cells.Item(75,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(75,2).Value = "0x01,0x90"
$ws.Cells.Item(75,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(75,4).Value = "0x01,0x90,"
$ws.Cells.Item(75,5).Value = "0x19"
$ws.Cells.Item(75,6).Value = 400
$ws.Cells.Item(75,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(75,8).Value = 400
$ws.Cells.Item(75,9).Value = 25

$ws.Cells.Item(76,1).Value = 45729.482181875
$ws.Cells.Item(76,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(76,2).Value = "0x01,0x90"
$ws.Cells.Item(76,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(76,4).Value = "0x01,0x90,"
$ws.Cells.Item(76,5).Value = "0x19"
$ws.Cells.Item(76,6).Value = 400
$ws.Cells.Item(76,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(76,8).Value = 400
$ws.Cells.Item(76,9).Value = 25

$ws.Cells.Item(77,1).Value = 45729.56561206019
$ws.Cells.Item(77,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(77,2).Value = "0x01,0x90"
$ws.Cells.Item(77,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(77,4).Value = "0x01,0x90,"
$ws.Cells.Item(77,5).Value = "0x19"
$ws.Cells.Item(77,6).Value = 400
$ws.Cells.Item(77,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(77,8).Value = 400
$ws.Cells.Item(77,9).Value = 25

$ws.Cells.Item(78,1).Value = 45729.56563403935
$ws.Cells.Item(78,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(78,2).Value = "0x01,0x90"
$ws.Cells.Item(78,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(78,4).Value = "0x01,0x90,"
$ws.Cells.Item(78,5).Value = "0x19"
$ws.Cells.Item(78,6).Value = 400
$ws.Cells.Item(78,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(78,8).Value = 400
$ws.Cells.Item(78,9).Value = 25

$ws.Cells.Item(79,1).Value = 45729.56565724537
$ws.Cells.Item(79,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(79,2).Value = "0x01,0x90"
$ws.Cells.Item(79,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(79,4).Value = "0x01,0x90,"
$ws.Cells.Item(79,5).Value = "0x19"
$ws.Cells.Item(79,6).Value = 400
$ws.Cells.Item(79,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(79,8).Value = 400
$ws.Cells.Item(79,9).Value = 25

$ws.Cells.Item(80,1).Value = 45729.64909206019
$ws.Cells.Item(80,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(80,2).Value = "0x01,0x90"
$ws.Cells.Item(80,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(80,4).Value = "0x01,0x90,"
$ws.Cells.Item(80,5).Value = "0x19"
$ws.Cells.Item(80,6).Value = 400
$ws.Cells.Item(80,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(80,8).Value = 400
$ws.Cells.Item(80,9).Value = 25

$ws.Cells.Item(81,1).Value = 45729.64911011574
$ws.Cells.Item(81,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(81,2).Value = "0x01,0x90"
$ws.Cells.Item(81,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(81,4).Value = "0x01,0x90,"
$ws.Cells.Item(81,5).Value = "0x19"
$ws.Cells.Item(81,6).Value = 400
$ws.Cells.Item(81,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(81,8).Value = 400
$ws.Cells.Item(81,9).Value = 25

$ws.Cells.Item(82,1).Value = 45729.64913337963
$ws.Cells.Item(82,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(82,2).Value = "0x01,0x90"
$ws.Cells.Item(82,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(82,4).Value = "0x01,0x90,"
$ws.Cells.Item(82,5).Value = "0x19"
$ws.Cells.Item(82,6).Value = 400
$ws.Cells.Item(82,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(82,8).Value = 400
$ws.Cells.Item(82,9).Value = 25

$ws.Cells.Item(83,1).Value = 45729.73256368055
$ws.Cells.Item(83,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(83,2).Value = "0x01,0x90"
$ws.Cells.Item(83,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(83,4).Value = "0x01,0x90,"
$ws.Cells.Item(83,5).Value = "0x19"
$ws.Cells.Item(83,6).Value = 400
$ws.Cells.Item(83,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(83,8).Value = 400
$ws.Cells.Item(83,9).Value = 25

$ws.Cells.Item(84,1).Value = 45729.73258552083
$ws.Cells.Item(84,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(84,2).Value = "0x01,0x90"
$ws.Cells.Item(84,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(84,4).Value = "0x01,0x90,"
$ws.Cells.Item(84,5).Value = "0x19"
$ws.Cells.Item(84,6).Value = 400
$ws.Cells.Item(84,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(84,8).Value = 400
$ws.Cells.Item(84,9).Value = 25

$ws.Cells.Item(85,1).Value = 45729.73260887731
$ws.Cells.Item(85,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(85,2).Value = "0x01,0x90"
$ws.Cells.Item(85,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(85,4).Value = "0x01,0x90,"
$ws.Cells.Item(85,5).Value = "0x19"
$ws.Cells.Item(85,6).Value = 400
$ws.Cells.Item(85,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(85,8).Value = 400
$ws.Cells.Item(85,9).Value = 25

# Sheet 3: ROW11-FE-LIFTER (rows 29-31)
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(29,1).Value = 45729.58035023148
$ws.Cells.Item(29,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(29,2).Value = "0x01,0x90"
$ws.Cells.Item(29,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
$ws.Cells.Item(29,4).Value = "0x01,0x90,"
$ws.Cells.Item(29,5).Value = "0x14"
$ws.Cells.Item(29,6).Value = 400
$ws.Cells.Item(29,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(29,8).Value = 400
$ws.Cells.Item(29,9).Value = 20

$ws.Cells.Item(30,1).Value = 45729.5803721875
$ws.Cells.Item(30,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(30,2).Value = "0x01,0x90"
$ws.Cells.Item(30,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
$ws.Cells.Item(30,4).Value = "0x01,0x90,"
$ws.Cells.Item(30,5).Value = "0x14"
$ws.Cells.Item(30,6).Value = 400
$ws.Cells.Item(30,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(30,8).Value = 400
$ws.Cells.Item(30,9).Value = 20

$ws.Cells.Item(31,1).Value = 45729.58039546297
$ws.Cells.Item(31,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(31,2).Value = "0x01,0x90"
$ws.Cells.Item(31,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
$ws.Cells.Item(31,4).Value = "0x01,0x90,"
$ws.Cells.Item(31,5).Value = "0x14"
$ws.Cells.Item(31,6).Value = 400
$ws.Cells.Item(31,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(31,8).Value = 400
$ws.Cells.Item(31,9).Value = 20

# Sheet 4: ROW11-MID-LIFTER (rows 68-85)
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(68,1).Value = 45729.31518523148
$ws.Cells.Item(68,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(68,2).Value = "0x01,0x90"
$ws.Cells.Item(68,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(68,4).Value = "0x01,0x90,"
$ws.Cells.Item(68,5).Value = "0x19"
$ws.Cells.Item(68,6).Value = 400
$ws.Cells.Item(68,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(68,8).Value = 400
$ws.Cells.Item(68,9).Value = 25

$ws.Cells.Item(69,1).Value = 45729.31520722222
$ws.Cells.Item(69,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(69,2).Value = "0x01,0x90"
$ws.Cells.Item(69,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(69,4).Value = "0x01,0x90,"
$ws.Cells.Item(69,5).Value = "0x19"
$ws.Cells.Item(69,6).Value = 400
$ws.Cells.Item(69,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(69,8).Value = 400
$ws.Cells.Item(69,9).Value = 25

$ws.Cells.Item(70,1).Value = 45729.3152303588
$ws.Cells.Item(70,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(70,2).Value = "0x01,0x90"
$ws.Cells.Item(70,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(70,4).Value = "0x01,0x90,"
$ws.Cells.Item(70,5).Value = "0x19"
$ws.Cells.Item(70,6).Value = 400
$ws.Cells.Item(70,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(70,8).Value = 400
$ws.Cells.Item(70,9).Value = 25

$ws.Cells.Item(71,1).Value = 45729.39866047454
$ws.Cells.Item(71,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(71,2).Value = "0x01,0x90"
$ws.Cells.Item(71,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(71,4).Value = "0x01,0x90,"
$ws.Cells.Item(71,5).Value = "0x19"
$ws.Cells.Item(71,6).Value = 400
$ws.Cells.Item(71,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(71,8).Value = 400
$ws.Cells.Item(71,9).Value = 25

$ws.Cells.Item(72,1).Value = 45729.39868246527
$ws.Cells.Item(72,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(72,2).Value = "0x01,0x90"
$ws.Cells.Item(72,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(72,4).Value = "0x01,0x90,"
$ws.Cells.Item(72,5).Value = "0x19"
$ws.Cells.Item(72,6).Value = 400
$ws.Cells.Item(72,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(72,8).Value = 400
$ws.Cells.Item(72,9).Value = 25

$ws.Cells.Item(73,1).Value = 45729.39870561343
$ws.Cells.Item(73,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(73,2).Value = "0x01,0x90"
$ws.Cells.Item(73,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(73,4).Value = "0x01,0x90,"
$ws.Cells.Item(73,5).Value = "0x19"
$ws.Cells.Item(73,6).Value = 400
$ws.Cells.Item(73,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(73,8).Value = 400
$ws.Cells.Item(73,9).Value = 25

$ws.Cells.Item(74,1).Value = 45729.48213673611
$ws.Cells.Item(74,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(74,2).Value = "0x01,0x90"
$ws.Cells.Item(74,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(74,4).Value = "0x01,0x90,"
$ws.Cells.Item(74,5).Value = "0x19"
$ws.Cells.Item(74,6).Value = 400
$ws.Cells.Item(74,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(74,8).Value = 400
$ws.Cells.Item(74,9).Value = 25

$ws.Cells.Item(75,1).Value = 45729.48215891204
$ws.Cells.Item(75,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(75,2).Value = "0x01,0x90"
$ws.Cells.Item(75,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(75,4).Value = "0x01,0x90,"
$ws.Cells.Item(75,5).Value = "0x19"
$ws.Cells.Item(75,6).Value = 400
$ws.Cells.Item(75,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(75,8).Value = 400
$ws.Cells.Item(75,9).Value = 25

$ws.Cells.Item(76,1).Value = 45729.482181875
$ws.Cells.Item(76,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(76,2).Value = "0x01,0x90"
$ws.Cells.Item(76,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(76,4).Value = "0x01,0x90,"
$ws.Cells.Item(76,5).Value = "0x19"
$ws.Cells.Item(76,6).Value = 400
$ws.Cells.Item(76,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(76,8).Value = 400
$ws.Cells.Item(76,9).Value = 25

$ws.Cells.Item(77,1).Value = 45729.56561206019
$ws.Cells.Item(77,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(77,2).Value = "0x01,0x90"
$ws.Cells.Item(77,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(77,4).Value = "0x01,0x90,"
$ws.Cells.Item(77,5).Value = "0x19"
$ws.Cells.Item(77,6).Value = 400
$ws.Cells.Item(77,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(77,8).Value = 400
$ws.Cells.Item(77,9).Value = 25

$ws.Cells.Item(78,1).Value = 45729.56563403935
$ws.Cells.Item(78,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(78,2).Value = "0x01,0x90"
$ws.Cells.Item(78,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(78,4).Value = "0x01,0x90,"
$ws.Cells.Item(78,5).Value = "0x19"
$ws.Cells.Item(78,6).Value = 400
$ws.Cells.Item(78,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(78,8).Value = 400
$ws.Cells.Item(78,9).Value = 25

$ws.Cells.Item(79,1).Value = 45729.56565724537
$ws.Cells.Item(79,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(79,2).Value = "0x01,0x90"
$ws.Cells.Item(79,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(79,4).Value = "0x01,0x90,"
$ws.Cells.Item(79,5).Value = "0x19"
$ws.Cells.Item(79,6).Value = 400
$ws.Cells.Item(79,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(79,8).Value = 400
$ws.Cells.Item(79,9).Value = 25

$ws.Cells.Item(80,1).Value = 45729.64909206019
$ws.Cells.Item(80,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(80,2).Value = "0x01,0x90"
$ws.Cells.Item(80,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(80,4).Value = "0x01,0x90,"
$ws.Cells.Item(80,5).Value = "0x19"
$ws.Cells.Item(80,6).Value = 400
$ws.Cells.Item(80,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(80,8).Value = 400
$ws.Cells.Item(80,9).Value = 25

$ws.Cells.Item(81,1).Value = 45729.64911011574
$ws.Cells.Item(81,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(81,2).Value = "0x01,0x90"
$ws.Cells.Item(81,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(81,4).Value = "0x01,0x90,"
$ws.Cells.Item(81,5).Value = "0x19"
$ws.Cells.Item(81,6).Value = 400
$ws.Cells.Item(81,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(81,8).Value = 400
$ws.Cells.Item(81,9).Value = 25

$ws.Cells.Item(82,1).Value = 45729.64913337963
$ws.Cells.Item(82,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(82,2).Value = "0x01,0x90"
$ws.Cells.Item(82,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(82,4).Value = "0x01,0x90,"
$ws.Cells.Item(82,5).Value = "0x19"
$ws.Cells.Item(82,6).Value = 400
$ws.Cells.Item(82,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(82,8).Value = 400
$ws.Cells.Item(82,9).Value = 25

$ws.Cells.Item(83,1).Value = 45729.73256368055
$ws.Cells.Item(83,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(83,2).Value = "0x01,0x90"
$ws.Cells.Item(83,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(83,4).Value = "0x01,0x90,"
$ws.Cells.Item(83,5).Value = "0x19"
$ws.Cells.Item(83,6).Value = 400
$ws.Cells.Item(83,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(83,8).Value = 400
$ws.Cells.Item(83,9).Value = 25

$ws.Cells.Item(84,1).Value = 45729.73258552083
$ws.Cells.Item(84,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(84,2).Value = "0x01,0x90"
$ws.Cells.Item(84,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(84,4).Value = "0x01,0x90,"
$ws.Cells.Item(84,5).Value = "0x19"
$ws.Cells.Item(84,6).Value = 400
$ws.Cells.Item(84,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(84,8).Value = 400
$ws.Cells.Item(84,9).Value = 25

$ws.Cells.Item(85,1).Value = 45729.73260887731
$ws.Cells.Item(85,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(85,2).Value = "0x01,0x90"
$ws.Cells.Item(85,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(85,4).Value = "0x01,0x90,"
$ws.Cells.Item(85,5).Value = "0x19"
$ws.Cells.Item(85,6).Value = 400
$ws.Cells.Item(85,7).Value = [double]"5.686312626471138e+23"
$ws.Cells.Item(85,8).Value = 400
$ws.Cells.Item(85,9).Value = 25

